# Insert a new data row at row 113 (shifting existing rows 113-168 down to
# 114-169), then populate it with the new "Friar" / bins (450 kilos) entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(113).Insert()

$ws.Range("A113").Value = 9
$ws.Range("B113").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C113").Value = "Metropolitana"
$ws.Range("D113").Value = 44978
$ws.Range("E113").Value = 13
$ws.Range("F113").Value = "Fruta"
$ws.Range("G113").Value = 100103
$ws.Range("H113").Value = "Frutos de hueso (carozo)"
$ws.Range("I113").Value = 100103002
$ws.Range("J113").Value = "Ciruela"
$ws.Range("K113").Value = "Friar"
$ws.Range("L113").Value = "Primera"
$ws.Range("M113").Value = 10
$ws.Range("N113").Value = 180000
$ws.Range("O113").Value = 180000
$ws.Range("P113").Value = 180000
$ws.Range("Q113").Value = "`$/bins (450 kilos)"
$ws.Range("R113").Value = "Región de O'Higgins"
$ws.Range("S113").Value = 400
$ws.Range("T113").Value = 450
